$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3847.1191
$ws.Range("I64").Value = 4134.9644
$ws.Range("K64").Value = 4134.9644
$ws.Range("M64").Value = -3886.9644
$ws.Range("H67").Value = 3847.1191
$ws.Range("I67").Value = 4134.9644
$ws.Range("K67").Value = 4134.9644
$ws.Range("M67").Value = -3276.9644
$ws.Range("H70").Value = 2744.875
$ws.Range("I70").Value = 3020.4
$ws.Range("J70").Value = 2619.6365
$ws.Range("K70").Value = 9061.200000000001
$ws.Range("L70").Value = 7858.9095
$ws.Range("M70").Value = -8791.200000000001
$ws.Range("N70").Value = -8398.9095
$ws.Range("H73").Value = 2744.875
$ws.Range("I73").Value = 3020.4
$ws.Range("J73").Value = 2619.6365
$ws.Range("K73").Value = 9061.200000000001
$ws.Range("L73").Value = 7858.9095
$ws.Range("M73").Value = -8125.200000000001
$ws.Range("N73").Value = -9730.9095
$ws.Range("H100").Value = 13890408
$ws.Range("I100").Value = 18519672
$ws.Range("J100").Value = 2618.6667
$ws.Range("K100").Value = 18519672
$ws.Range("L100").Value = 2618.6667
$ws.Range("M100").Value = -18519131
$ws.Range("N100").Value = -3700.6667
$ws.Range("H132").Value = 1132.5098
$ws.Range("I132").Value = 940.29785
$ws.Range("K132").Value = 2820.89355
$ws.Range("M132").Value = -290.8935500000002
$ws.Range("H138").Value = 2329.8572
$ws.Range("J138").Value = 2764.5925
$ws.Range("L138").Value = 8293.7775
$ws.Range("N138").Value = -18573.7775
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7818.329
$ws.Range("I32").Value = 5614.7573
$ws.Range("J32").Value = 20672.5
$ws.Range("K32").Value = 5614.7573
$ws.Range("L32").Value = 20672.5
$ws.Range("M32").Value = -5327.7573
$ws.Range("N32").Value = -21246.5
$ws.Range("H61").Value = 4582.302
$ws.Range("I61").Value = 6146.75
$ws.Range("J61").Value = 3287.5862
$ws.Range("K61").Value = 6146.75
$ws.Range("L61").Value = 3287.5862
$ws.Range("M61").Value = -5934.75
$ws.Range("N61").Value = -3711.5862
$ws.Range("H122").Value = 1285661
$ws.Range("I122").Value = 1976632.2
$ws.Range("J122").Value = 2428.5715
$ws.Range("K122").Value = 5929896.6
$ws.Range("L122").Value = 7285.7145
$ws.Range("M122").Value = -5927446.6
$ws.Range("N122").Value = -12185.7145
$ws.Range("H136").Value = 4582.302
$ws.Range("I136").Value = 6146.75
$ws.Range("J136").Value = 3287.5862
$ws.Range("K136").Value = 18440.25
$ws.Range("L136").Value = 9862.758600000001
$ws.Range("M136").Value = -15890.25
$ws.Range("N136").Value = -14962.7586
$ws.Range("H137").Value = 40780
$ws.Range("J137").Value = 40780
$ws.Range("L137").Value = 40780
$ws.Range("N137").Value = -50980
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 44279
$ws.Range("J53").Value = 44279
$ws.Range("L53").Value = 44279
$ws.Range("N53").Value = -45427
$ws.Range("H58").Value = 11450
$ws.Range("J58").Value = 11450
$ws.Range("L58").Value = 11450
$ws.Range("N58").Value = -12038
$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -46694
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("N52").Value = -50588
$ws.Range("H99").Value = 6258887
$ws.Range("I99").Value = 20310.125
$ws.Range("J99").Value = 10417938
$ws.Range("K99").Value = 20310.125
$ws.Range("L99").Value = 10417938
$ws.Range("M99").Value = -18812.125
$ws.Range("N99").Value = -10420934
$ws.Range("H105").Value = 2553.4348
$ws.Range("I105").Value = 2428.95
$ws.Range("K105").Value = 2428.95
$ws.Range("M105").Value = -681.9499999999998
$ws.Range("H122").Value = 1029971
$ws.Range("I122").Value = 1544134.1
$ws.Range("K122").Value = 4632402.300000001
$ws.Range("M122").Value = -4629952.300000001
$ws.Range("H126").Value = 6258887
$ws.Range("I126").Value = 20310.125
$ws.Range("J126").Value = 10417938
$ws.Range("K126").Value = 60930.375
$ws.Range("L126").Value = 31253814
$ws.Range("M126").Value = -58460.375
$ws.Range("N126").Value = -31258754
$ws.Range("H137").Value = 44000
$ws.Range("J137").Value = 44000
$ws.Range("L137").Value = 44000
$ws.Range("N137").Value = -54200
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 71428920
$ws.Range("J26").Value = 100000440
$ws.Range("L26").Value = 300001320
$ws.Range("N26").Value = -300001896
$ws.Range("H109").Value = 1776.3529
$ws.Range("J109").Value = 2375
$ws.Range("L109").Value = 7125
$ws.Range("N109").Value = -9205
$ws.Range("H131").Value = 1334312.6
$ws.Range("I131").Value = 12500375
$ws.Range("K131").Value = 37501125
$ws.Range("M131").Value = -37496085
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1180.5186
$ws.Range("I97").Value = 1195.1538
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 1195.1538
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -699.1538
$ws.Range("N97").Value = -1792
$ws.Range("H123").Value = 18646.617
$ws.Range("J123").Value = 18938.94
$ws.Range("L123").Value = 18938.94
$ws.Range("N123").Value = -23838.94
$ws.Range("H132").Value = 3038.9546
$ws.Range("I132").Value = 3535.182
$ws.Range("J132").Value = 2542.7273
$ws.Range("K132").Value = 10605.546
$ws.Range("L132").Value = 7628.1819
$ws.Range("M132").Value = -8075.545999999998
$ws.Range("N132").Value = -12688.1819
$ws.Range("H137").Value = 45726.668
$ws.Range("J137").Value = 45726.668
$ws.Range("L137").Value = 45726.668
$ws.Range("N137").Value = -55926.668
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 24392042
$ws.Range("I40").Value = 29413558
$ws.Range("K40").Value = 29413558
$ws.Range("M40").Value = -29413422
$ws.Range("H132").Value = 18526446
$ws.Range("I132").Value = 25651118
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 76953354
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -76950824
$ws.Range("N132").Value = -11960
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2916.077
$ws.Range("I122").Value = 2317.3333
$ws.Range("J122").Value = 3429.2856
$ws.Range("K122").Value = 6951.999899999999
$ws.Range("L122").Value = 10287.8568
$ws.Range("M122").Value = -4501.999899999999
$ws.Range("N122").Value = -15187.8568
$ws.Range("H132").Value = 2498.5557
$ws.Range("I132").Value = 1748.4166
$ws.Range("J132").Value = 3998.8333
$ws.Range("K132").Value = 5245.2498
$ws.Range("L132").Value = 11996.4999
$ws.Range("M132").Value = -2715.2498
$ws.Range("N132").Value = -17056.4999
